{"js": "// The document's single table has five \"data\" rows (0-based indices\n// 0, 4, 9, 14, 19), each holding five three-digit-by-one-digit\n// multiplication problems (one per column). We overwrite each cell's\n// text in place, addressed by its (row, col) position, rather than by\n// searching for the old text - some of the new answers duplicate text\n// that already exists elsewhere in the table (e.g. \"979\u00d78=7832\" is the\n// old value in row 0 and also becomes the new value of row 4, col 4),\n// so a global text search-and-replace could clobber the wrong cell.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst updates = [\n  { row: 0, col: 0, text: \"112\u00d77=784\" },\n  { row: 0, col: 1, text: \"151\u00d79=1359\" },\n  { row: 0, col: 2, text: \"153\u00d77=1071\" },\n  { row: 0, col: 3, text: \"474\u00d77=3318\" },\n  { row: 0, col: 4, text: \"784\u00d77=5488\" },\n\n  { row: 4, col: 0, text: \"202\u00d72=404\" },\n  { row: 4, col: 1, text: \"105\u00d76=630\" },\n  { row: 4, col: 2, text: \"525\u00d76=3150\" },\n  { row: 4, col: 3, text: \"302\u00d75=1510\" },\n  { row: 4, col: 4, text: \"979\u00d78=7832\" },\n\n  { row: 9, col: 0, text: \"562\u00d79=5058\" },\n  { row: 9, col: 1, text: \"311\u00d76=1866\" },\n  { row: 9, col: 2, text: \"329\u00d79=2961\" },\n  { row: 9, col: 3, text: \"720\u00d75=3600\" },\n  { row: 9, col: 4, text: \"622\u00d79=5598\" },\n\n  { row: 14, col: 0, text: \"755\u00d76=4530\" },\n  { row: 14, col: 1, text: \"613\u00d74=2452\" },\n  { row: 14, col: 2, text: \"580\u00d78=4640\" },\n  { row: 14, col: 3, text: \"163\u00d76=978\" },\n  { row: 14, col: 4, text: \"166\u00d73=498\" },\n\n  { row: 19, col: 0, text: \"627\u00d78=5016\" },\n  { row: 19, col: 1, text: \"314\u00d76=1884\" },\n  { row: 19, col: 2, text: \"594\u00d73=1782\" },\n  { row: 19, col: 3, text: \"455\u00d77=3185\" },\n  { row: 19, col: 4, text: \"528\u00d77=3696\" },\n];\n\nfor (const { row, col, text } of updates) {\n  const cell = table.getCell(row, col);\n  cell.value = text;\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# The table has five \"data\" rows (1, 5, 10, 15, 20), each with five\n# columns, holding a three-digit-by-one-digit multiplication problem. The\n# new values below replace the old ones cell-by-cell (by table position),\n# which is important because some new values duplicate text that appears\n# elsewhere in the document (e.g. \"979x8=7832\" is both an old answer in row\n# 1 and the new answer for row 5 col 5) - a plain global text search/replace\n# would misfire on those, so we address each cell directly.\n\n$values = @(\n    @(1,  1, \"112\u00d77=784\"),\n    @(1,  2, \"151\u00d79=1359\"),\n    @(1,  3, \"153\u00d77=1071\"),\n    @(1,  4, \"474\u00d77=3318\"),\n    @(1,  5, \"784\u00d77=5488\"),\n\n    @(5,  1, \"202\u00d72=404\"),\n    @(5,  2, \"105\u00d76=630\"),\n    @(5,  3, \"525\u00d76=3150\"),\n    @(5,  4, \"302\u00d75=1510\"),\n    @(5,  5, \"979\u00d78=7832\"),\n\n    @(10, 1, \"562\u00d79=5058\"),\n    @(10, 2, \"311\u00d76=1866\"),\n    @(10, 3, \"329\u00d79=2961\"),\n    @(10, 4, \"720\u00d75=3600\"),\n    @(10, 5, \"622\u00d79=5598\"),\n\n    @(15, 1, \"755\u00d76=4530\"),\n    @(15, 2, \"613\u00d74=2452\"),\n    @(15, 3, \"580\u00d78=4640\"),\n    @(15, 4, \"163\u00d76=978\"),\n    @(15, 5, \"166\u00d73=498\"),\n\n    @(20, 1, \"627\u00d78=5016\"),\n    @(20, 2, \"314\u00d76=1884\"),\n    @(20, 3, \"594\u00d73=1782\"),\n    @(20, 4, \"455\u00d77=3185\"),\n    @(20, 5, \"528\u00d77=3696\")\n)\n\nforeach ($item in $values) {\n    $rowIdx = $item[0]\n    $colIdx = $item[1]\n    $newText = $item[2]\n    $cell = $t.Cell($rowIdx, $colIdx)\n    $cell.Range.Text = $newText\n}\n"}
